$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$lastRow = $ws.UsedRange.Rows.Count

for ($r = 2; $r -le $lastRow; $r++) {
    $statutCell = $ws.Cells.Item($r, 1)
    $statutVal = $statutCell.Value()
    if ($statutVal -eq "⬛") {
        $statutCell.Value = "📘"
    }

    $labelCell = $ws.Cells.Item($r, 2)
    $labelVal = $labelCell.Value()
    if ($labelVal -eq "noir") {
        $labelCell.Value = "bleu"
    }
}
